$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1489
$ws.Range("I2").Value = 168.25
$ws.Range("K2").Value = 168.25
$ws.Range("M2").Value = -55.25
$ws.Range("H18").Value = 799.5
$ws.Range("J18").Value = 1000
$ws.Range("L18").Value = 1000
$ws.Range("N18").Value = -1568
$ws.Range("H43").Value = 4200.4
$ws.Range("J43").Value = 4200.4
$ws.Range("L43").Value = 4200.4
$ws.Range("N43").Value = -4338.4
$ws.Range("H106").Value = 1888.6
$ws.Range("I106").Value = 5450
$ws.Range("J106").Value = 998.25
$ws.Range("K106").Value = 5450
$ws.Range("L106").Value = 998.25
$ws.Range("M106").Value = -4819
$ws.Range("N106").Value = -2260.25
$ws.Range("H112").Value = 5770.778
$ws.Range("J112").Value = 6202.125
$ws.Range("L112").Value = 18606.375
$ws.Range("N112").Value = -20822.375
$ws.Range("H115").Value = 851.6667
$ws.Range("I115").Value = 277.5
$ws.Range("J115").Value = 2000
$ws.Range("K115").Value = 832.5
$ws.Range("L115").Value = 6000
$ws.Range("M115").Value = 734.5
$ws.Range("N115").Value = -9134
$ws.Range("H116").Value = 4777.2856
$ws.Range("I116").Value = 4923.5
$ws.Range("J116").Value = 4582.3335
$ws.Range("K116").Value = 4923.5
$ws.Range("L116").Value = 4582.3335
$ws.Range("M116").Value = -1481.5
$ws.Range("N116").Value = -11466.3335
$ws.Range("H138").Value = 6185.7104
$ws.Range("I138").Value = 3077.318
$ws.Range("J138").Value = 10459.75
$ws.Range("K138").Value = 9231.954000000002
$ws.Range("L138").Value = 31379.25
$ws.Range("M138").Value = -4091.954000000002
$ws.Range("N138").Value = -41659.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3727.5862
$ws.Range("I32").Value = 3196.4912
$ws.Range("K32").Value = 3196.4912
$ws.Range("M32").Value = -2909.4912
$ws.Range("H45").Value = 6308.933
$ws.Range("I45").Value = 5615.1665
$ws.Range("K45").Value = 5615.1665
$ws.Range("M45").Value = -5238.1665
$ws.Range("H61").Value = 3168.7454
$ws.Range("I61").Value = 2188.7273
$ws.Range("K61").Value = 2188.7273
$ws.Range("M61").Value = -1976.7273
$ws.Range("H122").Value = 2573.913
$ws.Range("I122").Value = 1370.6774
$ws.Range("K122").Value = 4112.0322
$ws.Range("M122").Value = -1662.0322
$ws.Range("H132").Value = 3063.9424
$ws.Range("I132").Value = 1731.9722
$ws.Range("J132").Value = 6060.875
$ws.Range("K132").Value = 5195.9166
$ws.Range("L132").Value = 18182.625
$ws.Range("M132").Value = -2665.9166
$ws.Range("N132").Value = -23242.625
$ws.Range("H136").Value = 3168.7454
$ws.Range("I136").Value = 2188.7273
$ws.Range("K136").Value = 6566.1819
$ws.Range("M136").Value = -4016.1819

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 943.7
$ws.Range("I11").Value = 551
$ws.Range("J11").Value = 1205.5
$ws.Range("K11").Value = 551
$ws.Range("L11").Value = 1205.5
$ws.Range("M11").Value = -411
$ws.Range("N11").Value = -1485.5
$ws.Range("H26").Value = 18117.5
$ws.Range("I26").Value = 18117.5
$ws.Range("K26").Value = 18117.5
$ws.Range("M26").Value = -17825.5
$ws.Range("H97").Value = 17999.6
$ws.Range("I97").Value = 9999.75
$ws.Range("K97").Value = 9999.75
$ws.Range("M97").Value = -9008.75
$ws.Range("H110").Value = 77000
$ws.Range("J110").Value = 77000
$ws.Range("L110").Value = 77000
$ws.Range("N110").Value = -85180
$ws.Range("H112").Value = 76000
$ws.Range("J112").Value = 76000
$ws.Range("L112").Value = 76000
$ws.Range("N112").Value = -78954

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 67030.336
$ws.Range("I60").Value = 1093
$ws.Range("J60").Value = 99999
$ws.Range("K60").Value = 1093
$ws.Range("L60").Value = 99999
$ws.Range("M60").Value = -582
$ws.Range("N60").Value = -101021
$ws.Range("H68").Value = 114824.5
$ws.Range("J68").Value = 114824.5
$ws.Range("L68").Value = 114824.5
$ws.Range("N68").Value = -116322.5
$ws.Range("H71").Value = 114824.5
$ws.Range("J71").Value = 114824.5
$ws.Range("L71").Value = 344473.5
$ws.Range("N71").Value = -351961.5
$ws.Range("H74").Value = 89700
$ws.Range("J74").Value = 91814.28999999999
$ws.Range("L74").Value = 91814.28999999999
$ws.Range("N74").Value = -93562.28999999999
$ws.Range("H77").Value = 89700
$ws.Range("J77").Value = 91814.28999999999
$ws.Range("L77").Value = 275442.87
$ws.Range("N77").Value = -284178.87
$ws.Range("H80").Value = 36500
$ws.Range("J80").Value = 36500
$ws.Range("L80").Value = 36500
$ws.Range("N80").Value = -38746
$ws.Range("H83").Value = 36500
$ws.Range("J83").Value = 36500
$ws.Range("L83").Value = 109500
$ws.Range("N83").Value = -120732
$ws.Range("H99").Value = 8605.75
$ws.Range("J99").Value = 7325.7
$ws.Range("L99").Value = 7325.7
$ws.Range("N99").Value = -10321.7
$ws.Range("H126").Value = 8605.75
$ws.Range("J126").Value = 7325.7
$ws.Range("L126").Value = 21977.1
$ws.Range("N126").Value = -26917.1
$ws.Range("H132").Value = 4770.241
$ws.Range("I132").Value = 3270
$ws.Range("J132").Value = 6616.6924
$ws.Range("K132").Value = 9810
$ws.Range("L132").Value = 19850.0772
$ws.Range("M132").Value = -7280
$ws.Range("N132").Value = -24910.0772
$ws.Range("H134").Value = 188866.14
$ws.Range("I134").Value = 2540.2942
$ws.Range("J134").Value = 505620.1
$ws.Range("K134").Value = 7620.882599999999
$ws.Range("L134").Value = 1516860.3
$ws.Range("M134").Value = -5085.882599999999
$ws.Range("N134").Value = -1521930.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 271766.34
$ws.Range("J5").Value = 13799
$ws.Range("L5").Value = 41397
$ws.Range("N5").Value = -41621
$ws.Range("H12").Value = 275.58823
$ws.Range("I12").Value = 109.75
$ws.Range("J12").Value = 326.6154
$ws.Range("K12").Value = 329.25
$ws.Range("L12").Value = 979.8462000000001
$ws.Range("M12").Value = -156.25
$ws.Range("N12").Value = -1325.8462
$ws.Range("H135").Value = 271766.34
$ws.Range("J135").Value = 13799
$ws.Range("L135").Value = 124191
$ws.Range("N135").Value = -129261

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 26584.285
$ws.Range("J15").Value = 26584.285
$ws.Range("L15").Value = 26584.285
$ws.Range("N15").Value = -27160.285
$ws.Range("H81").Value = 26584.285
$ws.Range("J81").Value = 26584.285
$ws.Range("L81").Value = 26584.285
$ws.Range("N81").Value = -28580.285
$ws.Range("H84").Value = 26584.285
$ws.Range("J84").Value = 26584.285
$ws.Range("L84").Value = 79752.855
$ws.Range("N84").Value = -89736.855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1425.25
$ws.Range("I16").Value = 1178.8572
$ws.Range("J16").Value = 3150
$ws.Range("K16").Value = 1178.8572
$ws.Range("L16").Value = 3150
$ws.Range("M16").Value = -1008.8572
$ws.Range("N16").Value = -3490
$ws.Range("H136").Value = 4079.2173
$ws.Range("I136").Value = 3368.8572
$ws.Range("K136").Value = 10106.5716
$ws.Range("M136").Value = -7556.571599999999
